$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the main title (A1): "Top 5 Research Activity Type Awards"
# -> "Top 5 International Research Activity Type Awards"
$ws.Range("A1").Value = "Top 5 International Research Activity Type Awards"

# Update the hidden/description title (A5): "Top 5 Research Activity Type Awards Description"
# -> "Top 5 International Research Activity Type Awards Description"
$ws.Range("A5").Value = "Top 5 International Research Activity Type Awards Description"

# The intro paragraph (A7) is no longer part of the merged A7:D7 block.
$ws.Range("A7:D7").UnMerge()

# Reset alignment on the unmerged cells to match the new (unmerged) layout.
$ws.Range("A7").HorizontalAlignment = 1
$ws.Range("B7:D7").HorizontalAlignment = 1
$ws.Range("E7:H7").HorizontalAlignment = 1
